$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increase offshore wind demand values for rows 29-31 (1450 -> 2900)
$ws.Range("E29").Value = 2900
$ws.Range("E30").Value = 2900
$ws.Range("E31").Value = 2900

# Move active selection to E31
$ws.Range("E31").Select() | Out-Null
